# #50 fixed signIn signOut, not manage go forward and go back on browser yet,
# change Navbar name to fake name
#
# On the ER diagram slide:
#  - the entity table previously labeled "JOIN" is renamed to
#    "USER_JOIN_EVENT" (it is the join table between USER and EVENT).
#  - the "1" cardinality label on the EVENT -> CATEGORY relationship
#    connector is changed to "*".

$p = $ppt.ActivePresentation

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)

    # Find the table shape whose header cell reads "JOIN" (there are several
    # tables named "Table 4" on this slide, so disambiguate by content) and
    # rename it to "USER_JOIN_EVENT".
    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        $shp = $s.Shapes.Item($i)
        if ($shp.HasTable) {
            $headerCell = $shp.Table.Cell(1, 1).Shape.TextFrame.TextRange
            if ($headerCell.Text -eq "JOIN") {
                $headerCell.Text = "USER_JOIN_EVENT"
            }
        }
    }

    # The cardinality textbox next to the EVENT -> CATEGORY connector is
    # uniquely named "TextBox 140".
    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        $shp = $s.Shapes.Item($i)
        if ($shp.Name -eq "TextBox 140") {
            $shp.TextFrame.TextRange.Text = "*"
        }
    }
}
